# "Adicionados balanços concatenados em uma única planilha."
# Rows 64 and 79 hold balance-sheet lines that are now blanked out (all
# the period columns C and E:AO cleared to empty text) as part of
# concatenating the balance sheets into the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($rowNum in 64, 79) {
    $r1 = $ws.Range("C$rowNum")
    $r1.Value = "'"
    $r1.Style = "Normal"

    $r2 = $ws.Range("E$($rowNum):AO$rowNum")
    $r2.Value = "'"
    $r2.Style = "Normal"
}
